$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Activity" (sheet2.xml): add a "reward 3" (ID/NUM) column pair and a
# new section row for activity 100001 / chapter 1.
# ---------------------------------------------------------------------------
$act = $wb.Worksheets.Item("Activity")

# Insert two new columns before the old "start time" column (J:K) so the
# existing start/end/duration columns shift right from J/K/L to L/M/N.
$act.Range("J1:K1").EntireColumn.Insert()

# Insert a new blank row before row 4 - the old rows 4/5/6 shift down to
# 5/6/7.
$act.Rows(4).Insert()

# New header cells for the inserted columns.
$act.Range("J1").Value = "奖励3ID"
$act.Range("K1").Value = "奖励3NUM"

# Row 2 (100001 / chapter 1 / section 1): reward2 is dropped, reward3 added.
$act.Range("H2").ClearContents()
$act.Range("I2").ClearContents()
$act.Range("J2").Value = 2007
$act.Range("K2").Value = 10

# Row 3 (100001 / chapter 1 / section 2): reward1 dropped, reward2 id
# changes, reward3 added.
$act.Range("F3").ClearContents()
$act.Range("G3").ClearContents()
$act.Range("H3").Value = 2003
$act.Range("J3").Value = 2007
$act.Range("K3").Value = 5

# Row 4 (new): 100001 / chapter 1, no section - only the lead-in columns are
# populated, everything else stays blank.
$act.Range("A4").Value = 100001
$act.Range("B4").Value = 1
$act.Range("C4").Value = "签到活动章1"
$act.Range("D4:E4").Clear()

# Row 5 (was row 4): 100001 / chapter 2 / section 1 - add reward3.
$act.Range("J5").Value = 2002
$act.Range("K5").Value = 3

# Row 6 (was row 5): 100002 / chapter 1 / section 1 - add reward3.
$act.Range("J6").Value = 2002
$act.Range("K6").Value = 2

# Row 7 (was row 6): 100003 / chapter 1 / section 1 - add reward3.
$act.Range("J7").Value = 2002
$act.Range("K7").Value = 2

$act.Range("M10").Select()

# ---------------------------------------------------------------------------
# Sheet "Item" (sheet1.xml): just a view/selection change.
# ---------------------------------------------------------------------------
$item = $wb.Worksheets.Item("Item")
$item.Activate()
$item.Range("G1").Select()
$item.Application.ActiveWindow.ScrollColumn = 7
$item.Range("K10").Select()

$act.Activate()
